$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InventoryList")

# Row 8: new product entry "In0005" / "Tov5" / "Sup3" ...
$ws.Range("C8").Value = "In0005"
$ws.Range("D8").Value = "Tov5"
$ws.Range("E8").Value = "Sup3"
$ws.Range("F8").Value = 33
$ws.Range("I8").Value = 45
$ws.Range("J8").Value = 2
$ws.Range("M8").Value = "Group3"

# Row 9: new product entry "In0006" / "Tov6" / "Sup3" (reuses supplier) ...
$ws.Range("C9").Value = "In0006"
$ws.Range("D9").Value = "Tov6"
$ws.Range("E9").Value = "Sup3"
$ws.Range("F9").Value = 33
$ws.Range("I9").Value = 32
$ws.Range("J9").Value = 44
$ws.Range("M9").Value = "Group3"
